$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Status column: handback is now in sync with en-US source ---
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$zhcn.Range("C2").Value = $newStatus
$dede.Range("C2").Value = $newStatus

# --- Latest Handback DateTime: refresh to the new handback run timestamps ---
$zhcn.Range("K2").Value = "2016-08-25 22:48:41"
$dede.Range("K2").Value = "2016-08-25 22:48:48"

# --- Error Detail: handback version mismatch is resolved, clear the error ---
$zhcn.Range("P2").Value = ""
$dede.Range("P2").Value = ""

# --- Column widths: widen Status columns and narrow the now-empty Error Detail columns ---
# (ColumnWidth is stored by Excel on a quantized pixel grid, so these inputs are chosen
# to land the stored <col width> as close as possible to the target layout.)
$statusColWidth = 29.166666666666668
$errorColWidth  = 12.833333333333332

$overview.Columns.Item(5).ColumnWidth = $statusColWidth
$overview.Columns.Item(6).ColumnWidth = $statusColWidth

$zhcn.Columns.Item(3).ColumnWidth = $statusColWidth
$zhcn.Columns.Item(16).ColumnWidth = $errorColWidth

$dede.Columns.Item(3).ColumnWidth = $statusColWidth
$dede.Columns.Item(16).ColumnWidth = $errorColWidth
